# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Thu Jun  1 08:14:15 UTC 2023 with GitHub Actions".
#
# Numeric-looking price strings (e.g. "1.001") must be written as literal
# text (matching the original inlineStr cells), not auto-converted to
# Excel numbers. We force text via NumberFormat "@" before assignment and
# then clear the formatting again so no extra style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.836.02"
$ws.Range("E2").Value = "  -1.31%  "

# Row 3
$ws.Range("D3").Value = "1.856.27"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.99%  "

# Row 6
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5051"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.94%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3632"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.51%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07154"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8906"
$ws.Range("D10").ClearFormats()

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.66"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.46%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07467"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.71%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.850.54"
$ws.Range("E13").Value = "  -1.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.48%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.218"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008492"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.00%  "

# Row 20
$ws.Range("D20").Value = "26.876.41"
$ws.Range("E20").Value = "  -1.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.008"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.42%  "

# Row 22
$ws.Range("D22").Value = "2.091.96"
$ws.Range("E22").Value = "  -0.41%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.417"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.38"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.20%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.795"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.40%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.047"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.64%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.632"
$ws.Range("D30").ClearFormats()

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.650"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09223"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05074"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.70%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.981"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7433"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.44%  "

# Row 36
$ws.Range("E36").Value = "  -2.39%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.251"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.22%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.504"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.72%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.096"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01981"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5316"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.87%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.468"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.47%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.403"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1458"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.76%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.07%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4631"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.84%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.05"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.555"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.72"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.07%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05934"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.45%  "

